# The trainer for the first block of topics (rows 2-44, i.e. range C2:C44)
# is renamed from "Karun Karthik" to "BhanuTeja Reddy" (matching the name
# already used for the rest of the plan). Excel will drop the now-unused
# "Karun Karthik" shared string automatically on save.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("C2:C44")
$rng.Value = "BhanuTeja Reddy"

# Leave the same range selected, matching the state captured in the saved
# workbook.
$rng.Select()
